$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1336404.6
$ws.Range("I15").Value = 1336404.6
$ws.Range("K15").Value = 4009213.8
$ws.Range("M15").Value = -4009044.8

$ws.Range("H28").Value = 456.81818
$ws.Range("I28").Value = 307.05884
$ws.Range("J28").Value = 966
$ws.Range("K28").Value = 307.05884
$ws.Range("L28").Value = 966
$ws.Range("M28").Value = 177.94116
$ws.Range("N28").Value = -1936

$ws.Range("H107").Value = 510.81818
$ws.Range("I107").Value = 371.25
$ws.Range("J107").Value = 883
$ws.Range("K107").Value = 371.25
$ws.Range("L107").Value = 883
$ws.Range("M107").Value = 1548.75
$ws.Range("N107").Value = -4723

$ws.Range("H131").Value = 35827.742
$ws.Range("I131").Value = 50702
$ws.Range("J131").Value = 8783.637000000001
$ws.Range("K131").Value = 152106
$ws.Range("L131").Value = 26350.911
$ws.Range("M131").Value = -147066
$ws.Range("N131").Value = -36430.911

$ws.Range("H132").Value = 7048.039
$ws.Range("I132").Value = 7426.6665
$ws.Range("J132").Value = 6507.143
$ws.Range("K132").Value = 22279.9995
$ws.Range("L132").Value = 19521.429
$ws.Range("M132").Value = -19749.9995
$ws.Range("N132").Value = -24581.429

$ws.Range("H137").Value = 1903.0555
$ws.Range("I137").Value = 2292.28
$ws.Range("J137").Value = 1018.4545
$ws.Range("K137").Value = 6876.84
$ws.Range("L137").Value = 3055.3635
$ws.Range("M137").Value = -4326.84
$ws.Range("N137").Value = -8155.3635

$ws.Range("H138").Value = 1463.1538
$ws.Range("I138").Value = 1391.0834
$ws.Range("J138").Value = 1578.4667
$ws.Range("K138").Value = 4173.2502
$ws.Range("L138").Value = 4735.4001
$ws.Range("M138").Value = 966.7497999999996
$ws.Range("N138").Value = -15015.4001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11400.685
$ws.Range("I32").Value = 6455.012
$ws.Range("J32").Value = 51955.2
$ws.Range("K32").Value = 6455.012
$ws.Range("L32").Value = 51955.2
$ws.Range("M32").Value = -6168.012
$ws.Range("N32").Value = -52529.2

$ws.Range("H132").Value = 1468.84
$ws.Range("I132").Value = 1326.3529
$ws.Range("J132").Value = 1617.1428
$ws.Range("K132").Value = 3979.0587
$ws.Range("L132").Value = 4851.428400000001
$ws.Range("M132").Value = -1449.0587
$ws.Range("N132").Value = -9911.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2749.7837
$ws.Range("I7").Value = 14349.429
$ws.Range("J7").Value = 43.2
$ws.Range("K7").Value = 14349.429
$ws.Range("L7").Value = 43.2
$ws.Range("M7").Value = -14236.429
$ws.Range("N7").Value = -269.2

$ws.Range("H22").Value = 335.10715
$ws.Range("I22").Value = 286.72223
$ws.Range("J22").Value = 422.2
$ws.Range("K22").Value = 286.72223
$ws.Range("L22").Value = 422.2
$ws.Range("M22").Value = 63.27776999999998
$ws.Range("N22").Value = -1122.2

$ws.Range("H31").Value = 7248991.5
$ws.Range("I31").Value = 2046.2667
$ws.Range("J31").Value = 20837014
$ws.Range("K31").Value = 2046.2667
$ws.Range("L31").Value = 20837014
$ws.Range("M31").Value = -1751.2667
$ws.Range("N31").Value = -20837604

$ws.Range("H34").Value = 7248991.5
$ws.Range("I34").Value = 2046.2667
$ws.Range("J34").Value = 20837014
$ws.Range("K34").Value = 2046.2667
$ws.Range("L34").Value = 20837014
$ws.Range("M34").Value = -1844.2667
$ws.Range("N34").Value = -20837418

$ws.Range("H58").Value = 727447.7
$ws.Range("I58").Value = 1510.6666
$ws.Range("J58").Value = 1671165.9
$ws.Range("K58").Value = 1510.6666
$ws.Range("L58").Value = 1671165.9
$ws.Range("M58").Value = -1307.6666
$ws.Range("N58").Value = -1671571.9

$ws.Range("H136").Value = 727447.7
$ws.Range("I136").Value = 1510.6666
$ws.Range("J136").Value = 1671165.9
$ws.Range("K136").Value = 4531.9998
$ws.Range("L136").Value = 5013497.699999999
$ws.Range("M136").Value = -1981.9998
$ws.Range("N136").Value = -5018597.699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50.7
$ws.Range("I2").Value = 47.42857
$ws.Range("J2").Value = 52.46154
$ws.Range("K2").Value = 284.57142
$ws.Range("L2").Value = 314.76924
$ws.Range("M2").Value = -171.57142
$ws.Range("N2").Value = -540.76924

$ws.Range("H33").Value = 194.72223
$ws.Range("I33").Value = 85.36364
$ws.Range("J33").Value = 366.57144
$ws.Range("K33").Value = 512.18184
$ws.Range("L33").Value = 2199.42864
$ws.Range("M33").Value = -229.18184
$ws.Range("N33").Value = -2765.42864

$ws.Range("H68").Value = 686.1667
$ws.Range("I68").Value = 370
$ws.Range("J68").Value = 844.25
$ws.Range("K68").Value = 1110
$ws.Range("L68").Value = 2532.75
$ws.Range("M68").Value = -299
$ws.Range("N68").Value = -4154.75

$ws.Range("H71").Value = 686.1667
$ws.Range("I71").Value = 370
$ws.Range("J71").Value = 844.25
$ws.Range("K71").Value = 3330
$ws.Range("L71").Value = 7598.25
$ws.Range("M71").Value = 726
$ws.Range("N71").Value = -15710.25

$ws.Range("H86").Value = 613.1111
$ws.Range("I86").Value = 445.42856
$ws.Range("J86").Value = 1200
$ws.Range("K86").Value = 1336.28568
$ws.Range("L86").Value = 3600
$ws.Range("M86").Value = -150.28568
$ws.Range("N86").Value = -5972

$ws.Range("H89").Value = 613.1111
$ws.Range("I89").Value = 445.42856
$ws.Range("J89").Value = 1200
$ws.Range("K89").Value = 4008.85704
$ws.Range("L89").Value = 10800
$ws.Range("M89").Value = 1919.14296
$ws.Range("N89").Value = -22656

$ws.Range("H131").Value = 790
$ws.Range("I131").Value = 346.25
$ws.Range("J131").Value = 1233.75
$ws.Range("K131").Value = 1038.75
$ws.Range("L131").Value = 3701.25
$ws.Range("M131").Value = 4001.25
$ws.Range("N131").Value = -13781.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2234.3076
$ws.Range("J7").Value = 2201.5
$ws.Range("L7").Value = 2201.5
$ws.Range("N7").Value = -2425.5

$ws.Range("H40").Value = 2340.439
$ws.Range("I40").Value = 1950.9354
$ws.Range("J40").Value = 3547.9
$ws.Range("K40").Value = 1950.9354
$ws.Range("L40").Value = 3547.9
$ws.Range("M40").Value = -1814.9354
$ws.Range("N40").Value = -3819.9

$ws.Range("H46").Value = 1429981.9
$ws.Range("I46").Value = 1820.25
$ws.Range("J46").Value = 3334197.2
$ws.Range("K46").Value = 1820.25
$ws.Range("L46").Value = 3334197.2
$ws.Range("M46").Value = -1632.25
$ws.Range("N46").Value = -3334573.2

$ws.Range("H68").Value = 3183.3333
$ws.Range("I68").Value = 2550
$ws.Range("J68").Value = 3500
$ws.Range("K68").Value = 2550
$ws.Range("L68").Value = 3500
$ws.Range("M68").Value = -1801
$ws.Range("N68").Value = -4998

$ws.Range("H71").Value = 3183.3333
$ws.Range("I71").Value = 2550
$ws.Range("J71").Value = 3500
$ws.Range("K71").Value = 12750
$ws.Range("L71").Value = 17500
$ws.Range("M71").Value = -9006
$ws.Range("N71").Value = -24988

$ws.Range("H93").Value = 2073.1765
$ws.Range("I93").Value = 1828.75
$ws.Range("J93").Value = 2659.8
$ws.Range("K93").Value = 1828.75
$ws.Range("L93").Value = 2659.8
$ws.Range("M93").Value = -580.75
$ws.Range("N93").Value = -5155.8

$ws.Range("H98").Value = 16450.8
$ws.Range("J98").Value = 16450.8
$ws.Range("L98").Value = 16450.8
$ws.Range("N98").Value = -22440.8

$ws.Range("H101").Value = 16833.334
$ws.Range("J101").Value = 16833.334
$ws.Range("L101").Value = 16833.334
$ws.Range("N101").Value = -23323.334

$ws.Range("H105").Value = 50614
$ws.Range("J105").Value = 50614
$ws.Range("L105").Value = 50614
$ws.Range("N105").Value = -57602

$ws.Range("H126").Value = 2234.3076
$ws.Range("J126").Value = 2201.5
$ws.Range("L126").Value = 6604.5
$ws.Range("N126").Value = -11544.5

$ws.Range("H136").Value = 1245.225
$ws.Range("I136").Value = 775.25
$ws.Range("J136").Value = 3125.125
$ws.Range("K136").Value = 2325.75
$ws.Range("L136").Value = 9375.375
$ws.Range("M136").Value = 224.25
$ws.Range("N136").Value = -14475.375
